$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-02-16 09:42:58"
$wsZh.Range("G4").Value = "2016-02-16 09:43:53"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-02-16 09:43:12"
$wsDe.Range("G4").Value = "2016-02-16 09:44:21"
